$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '62.707.31'
$ws.Range('E2').Value = '  -0.99%  '

$ws.Range('D3').Value = '3.024.45'
$ws.Range('E3').Value = '  -1.22%  '

$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.18%  '

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '586.33'
$ws.Range('E5').Value = '  -0.25%  '

$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '147.72'
$ws.Range('E6').Value = '  -4.64%  '

$ws.Range('E7').Value = '  -0.08%  '

$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.523'
$ws.Range('E8').Value = '  -2.55%  '

$ws.Range('D9').Value = '3.022.94'
$ws.Range('E9').Value = '  -1.20%  '

$ws.Range('E10').Value = '  -3.36%  '

$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '5.80'
$ws.Range('E11').Value = '  -0.16%  '

$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.443'
$ws.Range('E12').Value = '  -1.16%  '

$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '0.0000231'
$ws.Range('E13').Value = '  -2.39%  '

$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '34.90'
$ws.Range('E14').Value = '  -5.10%  '

$ws.Range('E15').Value = '  +2.35%  '

$ws.Range('D16').Value = '3.516.14'
$ws.Range('E16').Value = '  -1.56%  '

$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '7.07'
$ws.Range('E17').Value = '  -0.60%  '

$ws.Range('D18').Value = '62.610.79'
$ws.Range('E18').Value = '  -1.21%  '

$ws.Range('D19').Value = '3.016.98'

$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '465.43'
$ws.Range('E20').Value = '  -0.92%  '

$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '13.98'
$ws.Range('E21').Value = '  -2.00%  '

$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '0.689'
$ws.Range('E22').Value = '  -1.98%  '

$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '7.38'
$ws.Range('E23').Value = '  -1.32%  '

$ws.Range('E24').Value = '  -1.73%  '

$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '80.11'
$ws.Range('E25').Value = '  -0.46%  '

$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '12.44'
$ws.Range('E26').Value = '  -2.44%  '

$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '10.19'
$ws.Range('E27').Value = '  -1.31%  '

$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '0.999'
$ws.Range('E28').Value = '  +0.12%  '

$ws.Range('E29').Value = '  -0.06%  '

$ws.Range('E30').Value = '  -0.61%  '

$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '7.20'
$ws.Range('E31').Value = '  -2.35%  '

$ws.Range('E32').Value = '  -0.58%  '

$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '27.50'
$ws.Range('E33').Value = '  +1.69%  '

$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.107'
$ws.Range('E34').Value = '  -3.80%  '

$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '1.04'
$ws.Range('E35').Value = '  -0.31%  '

$ws.Range('D36').Value = '0.0₃0798'
$ws.Range('E36').Value = '  -1.83%  '

$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '5.77'
$ws.Range('E37').Value = '  -3.19%  '

$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '2.14'
$ws.Range('E38').Value = '  -2.51%  '

$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '50.52'
$ws.Range('E39').Value = '  +0.10%  '

$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '9.01'
$ws.Range('E40').Value = '  -2.05%  '

$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '2.95'
$ws.Range('E41').Value = '  -8.86%  '

$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '423.98'
$ws.Range('E42').Value = '  -2.25%  '

$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.113'
$ws.Range('E43').Value = '  +1.07%  '

$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.278'
$ws.Range('E44').Value = '  -1.83%  '

$ws.Range('D45').Value = '2.788.06'
$ws.Range('E45').Value = '  -0.16%  '

$ws.Range('E46').Value = '  -1.30%  '

$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '38.18'
$ws.Range('E47').Value = '  -6.00%  '

$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '129.57'
$ws.Range('E48').Value = '  +0.13%  '

$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '24.18'
$ws.Range('E50').Value = '  -3.13%  '

$ws.Range('E51').Value = '  -0.59%  '
